# Weekly update: insert the latest week's Chirimoya price record at the
# top of the data block (row 4), pushing all existing data rows down by
# one (old row 4 -> row 5, old row 24 -> row 25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current first data row (row 4). Excel shifts
# rows 4:24 down to 5:25 and extends the used range accordingly.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest week's record.
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44530
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 2000
$ws.Range("O4").Value = 2100
$ws.Range("P4").Value = 2050
$ws.Range("Q4").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 2050
$ws.Range("T4").Value = 1
